# Update countries & provincias Spain
#
# The "Pais" sheet lists countries with daily COVID-19 stats
# (Casos totales, Nuevos casos, Casos activos, Recuperados,
# Casos criticos, Muertes hoy, Muertes) sorted descending by
# "Casos totales" (column B). This applies an updated day's data for
# a handful of countries. Because the sheet stays sorted by total
# cases, a couple of pairs of adjacent rows swap places once the new
# totals are plugged in (their labels/values trade rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Chequia (row 46) - refreshed Recuperados/Casos activos/Muertes hoy,
# no reordering needed here.
$ws.Range("D46").Value = 3471
$ws.Range("E46").Value = 4039
$ws.Range("F46").Value = 62

# Malasia now has more total cases (6298) than Egipto (6193), so it
# moves above Egipto: row 52 becomes Malasia's updated figures, and
# Egipto's (unchanged) figures shift down to row 53.
$ws.Range("A52").Value = "Malasia"
$ws.Range("B52").Value = 6298
$ws.Range("C52").Value = 122
$ws.Range("D52").Value = 4413
$ws.Range("E52").Value = 1780
$ws.Range("F52").Value = 27
$ws.Range("G52").Value = 2
$ws.Range("H52").Value = 105

$ws.Range("A53").Value = "Egipto"
$ws.Range("B53").Value = 6193
$ws.Range("C53").Value = 0
$ws.Range("D53").Value = 1522
$ws.Range("E53").Value = 4256
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 415

# Etiopia's updated total (135) now ties/exceeds Madagascar's, so it
# moves above Madagascar: row 144 becomes Etiopia's updated figures,
# and Madagascar's (unchanged) figures shift down to row 145.
$ws.Range("A144").Value = "Etiopia"
$ws.Range("B144").Value = 135
$ws.Range("C144").Value = 2
$ws.Range("D144").Value = 75
$ws.Range("E144").Value = 57
$ws.Range("F144").Value = 0
$ws.Range("G144").Value = 0
$ws.Range("H144").Value = 3

$ws.Range("A145").Value = "Madagascar"
$ws.Range("B145").Value = 135
$ws.Range("C145").Value = 0
$ws.Range("D145").Value = 97
$ws.Range("E145").Value = 38
$ws.Range("F145").Value = 1
$ws.Range("G145").Value = 0
$ws.Range("H145").Value = 0

# San Vicente y las Granadinas and Namibia are tied (16 total cases
# each); San Vicente now sorts first, so rows 194/195 swap labels
# (the underlying figures are identical for both, so only the
# country names trade places).
$ws.Range("A194").Value = "San Vicente y las Granadinas"
$ws.Range("B194").Value = 16
$ws.Range("C194").Value = 0
$ws.Range("D194").Value = 8
$ws.Range("E194").Value = 8
$ws.Range("F194").Value = 0
$ws.Range("G194").Value = 0
$ws.Range("H194").Value = 0

$ws.Range("A195").Value = "Namibia"
$ws.Range("B195").Value = 16
$ws.Range("C195").Value = 0
$ws.Range("D195").Value = 8
$ws.Range("E195").Value = 8
$ws.Range("F195").Value = 0
$ws.Range("G195").Value = 0
$ws.Range("H195").Value = 0
